# Update "想去人数" (number of people wanting to attend) figures on the
# "展览" and "全部类型" sheets, as produced by the gh-pages data refresh.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 481
    $ws.Range("F8").Value = 1348
    $ws.Range("F9").Value = 3987
    $ws.Range("F10").Value = 87
}
